$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 24.08000000000033
$ws.Range("H2").Value = 0.000009673536133014693
$ws.Range("I2").Value = 0.000009673536133014693
$ws.Range("L2").Value = 55.94790282046719
$ws.Range("M2").Value = "[29.450986528447928, 82.44481911248646]"
$ws.Range("N2").Value = 0.0001052924133282396
$ws.Range("O2").Value = 0.0001052924133282396
$ws.Range("P2").Value = 1.226447582482502
$ws.Range("Q2").Value = "[0.6855527512338089, 1.7673424137311944]"
$ws.Range("R2").Value = 0.00003831238604945497
$ws.Range("S2").Value = 0.00003831238604945497
$ws.Range("T2").Value = 65.72112959929896
$ws.Range("U2").Value = "[51.205237370971375, 80.23702182762655]"
$ws.Range("V2").Value = 0.00000000000861466453727644
$ws.Range("W2").Value = 0.00000000000861466453727644
$ws.Range("X2").Value = 19.37969969969996
$ws.Range("Y2").Value = 17.30674674674698
$ws.Range("Z2").Value = 21.45265265265294

# Row 3
$ws.Range("F3").Value = 24.08000000000033
$ws.Range("H3").Value = 0.00001226671726162198
$ws.Range("I3").Value = 0.00001226671726162198
$ws.Range("L3").Value = 51.5709127334451
$ws.Range("M3").Value = "[26.405830104188297, 76.7359953627019]"
$ws.Range("N3").Value = 0.0001565063980508707
$ws.Range("O3").Value = 0.0001565063980508707
$ws.Range("P3").Value = 2.081816152829272
$ws.Range("Q3").Value = "[1.5786581702723472, 2.5849741353861964]"
$ws.Range("R3").Value = 0.0000000001133624305538206
$ws.Range("S3").Value = 0.0000000001133624305538206
$ws.Range("T3").Value = 60.96993183011804
$ws.Range("U3").Value = "[47.40922339720848, 74.5306402630276]"
$ws.Range("V3").Value = 0.00000000001057887111244327
$ws.Range("W3").Value = 0.00000000001057887111244327
$ws.Range("X3").Value = 16.10154154154176
$ws.Range("Y3").Value = 14.1732132132134
$ws.Range("Z3").Value = 18.02986986987012

# Row 4
$ws.Range("F4").Value = 24.08000000000033
$ws.Range("H4").Value = 0.00008527500181454428
$ws.Range("I4").Value = 0.00008527500181454428
$ws.Range("L4").Value = 63.95230906210087
$ws.Range("M4").Value = "[34.67125996598391, 93.23335815821783]"
$ws.Range("N4").Value = 0.00006595440847823042
$ws.Range("O4").Value = 0.00006595440847823042
$ws.Range("P4").Value = 1.591237119836271
$ws.Range("Q4").Value = "[0.9748685912040393, 2.207605648468503]"
$ws.Range("R4").Value = 0.000004718438163964223
$ws.Range("S4").Value = 0.000004718438163964223
$ws.Range("T4").Value = 68.97083408855211
$ws.Range("U4").Value = "[50.06347676613868, 87.87819141096554]"
$ws.Range("V4").Value = 0.00000000312350123543581
$ws.Range("W4").Value = 0.00000000312350123543581
$ws.Range("X4").Value = 17.98166166166191
$ws.Range("Y4").Value = 15.61945945945967
$ws.Range("Z4").Value = 20.34386386386414

# Row 5
$ws.Range("F5").Value = 24.08000000000033
$ws.Range("H5").Value = 0.0008543181850799719
$ws.Range("I5").Value = 0.0008543181850799719
$ws.Range("L5").Value = 42.86916739963277
$ws.Range("M5").Value = "[16.751833381822962, 68.98650141744257]"
$ws.Range("N5").Value = 0.001864770619364631
$ws.Range("O5").Value = 0.001864770619364631
$ws.Range("P5").Value = 1.364816027685655
$ws.Range("Q5").Value = "[0.635236952978115, 2.094395102393195]"
$ws.Range("R5").Value = 0.0004765024184938405
$ws.Range("S5").Value = 0.0004765024184938405
$ws.Range("T5").Value = 61.81000162425103
$ws.Range("U5").Value = "[46.81937941053786, 76.8006238379642]"
$ws.Range("V5").Value = 0.00000000012470802168707
$ws.Range("W5").Value = 0.00000000012470802168707
$ws.Range("X5").Value = 18.84940940940967
$ws.Range("Y5").Value = 16.05333333333355
$ws.Range("Z5").Value = 21.64548548548578

# Row 6
$ws.Range("F6").Value = 24.08000000000033
$ws.Range("H6").Value = 0.0009592802262278211
$ws.Range("I6").Value = 0.0009592802262278211
$ws.Range("L6").Value = 38.10735167269533
$ws.Range("M6").Value = "[15.197442507176156, 61.017260838214504]"
$ws.Range("N6").Value = 0.001642285833519663
$ws.Range("O6").Value = 0.001642285833519663
$ws.Range("P6").Value = 1.943447707626119
$ws.Range("Q6").Value = "[1.2516054816103486, 2.6352899336418885]"
$ws.Range("R6").Value = 0.000001002907303648826
$ws.Range("S6").Value = 0.000001002907303648826
$ws.Range("T6").Value = 52.98378574584978
$ws.Range("U6").Value = "[39.5355680524981, 66.43200343920145]"
$ws.Range("V6").Value = 0.0000000004283216004097312
$ws.Range("W6").Value = 0.0000000004283216004097312
$ws.Range("X6").Value = 16.63183183183206
$ws.Range("Y6").Value = 13.98038038038057
$ws.Range("Z6").Value = 19.28328328328354

# Row 7
$ws.Range("F7").Value = 24.08000000000033
$ws.Range("H7").Value = 0.0009612269819817554
$ws.Range("I7").Value = 0.0009612269819817554
$ws.Range("L7").Value = 44.36714567818117
$ws.Range("M7").Value = "[17.58338481338663, 71.1509065429757]"
$ws.Range("N7").Value = 0.001708979690091539
$ws.Range("O7").Value = 0.001708979690091539
$ws.Range("P7").Value = 1.415131825941349
$ws.Range("Q7").Value = "[0.6855527512338098, 2.144710900648888]"
$ws.Range("R7").Value = 0.0003114669904349388
$ws.Range("S7").Value = 0.0003114669904349388
$ws.Range("T7").Value = 61.19996028937858
$ws.Range("U7").Value = "[45.52093414361265, 76.87898643514451]"
$ws.Range("V7").Value = 0.0000000005484375176223466
$ws.Range("W7").Value = 0.0000000005484375176223466
$ws.Range("X7").Value = 18.65657657657682
$ws.Range("Y7").Value = 15.86050050050071
$ws.Range("Z7").Value = 21.45265265265294

# Row 8
$ws.Range("F8").Value = 24.08000000000033
$ws.Range("H8").Value = 0.00000112908819671631
$ws.Range("I8").Value = 0.00000112908819671631
$ws.Range("L8").Value = 54.20521887674822
$ws.Range("M8").Value = "[34.14516612698789, 74.26527162650855]"
$ws.Range("N8").Value = 0.000002082265026182029
$ws.Range("O8").Value = 0.000002082265026182029
$ws.Range("P8").Value = 1.427710775505271
$ws.Range("Q8").Value = "[0.9748685912040393, 1.8805529598065025]"
$ws.Range("R8").Value = 0.00000009423575653144667
$ws.Range("S8").Value = 0.00000009423575653144667
$ws.Range("T8").Value = 52.14334089256764
$ws.Range("U8").Value = "[39.55351404664714, 64.73316773848813]"
$ws.Range("V8").Value = 0.0000000001102196112157117
$ws.Range("W8").Value = 0.0000000001102196112157117
$ws.Range("X8").Value = 18.60836836836862
$ws.Range("Y8").Value = 16.8728728728731
$ws.Range("Z8").Value = 20.34386386386414

# Row 9
$ws.Range("F9").Value = 24.08000000000033
$ws.Range("H9").Value = 0.000002235891482960106
$ws.Range("I9").Value = 0.000002235891482960106
$ws.Range("L9").Value = 52.49504470195755
$ws.Range("M9").Value = "[31.78795459055918, 73.20213481335593]"
$ws.Range("N9").Value = 0.000006458154654387371
$ws.Range("O9").Value = 0.000006458154654387371
$ws.Range("P9").Value = 1.478026573760964
$ws.Range("Q9").Value = "[1.0126054398958084, 1.9434477076261194]"
$ws.Range("R9").Value = 0.00000008046485588941721
$ws.Range("S9").Value = 0.00000008046485588941721
$ws.Range("T9").Value = 60.82147984980447
$ws.Range("U9").Value = "[48.22642361452404, 73.4165360850849]"
$ws.Range("V9").Value = 0.000000000001232125512728999
$ws.Range("W9").Value = 0.000000000001232125512728999
$ws.Range("X9").Value = 18.41553553553578
$ws.Range("Y9").Value = 16.63183183183205
$ws.Range("Z9").Value = 20.19923923923951

# Row 10
$ws.Range("F10").Value = 23.42000000000022
$ws.Range("H10").Value = 0.0000004490228553422071
$ws.Range("I10").Value = 0.0000004490228553422071
$ws.Range("L10").Value = 64.07218497607408
$ws.Range("M10").Value = "[41.085122969841876, 87.05924698230629]"
$ws.Range("N10").Value = 0.000001164107659379354
$ws.Range("O10").Value = 0.000001164107659379354
$ws.Range("P10").Value = 1.465447624197041
$ws.Range("Q10").Value = "[1.0377633390236554, 1.8931319093704264]"
$ws.Range("R10").Value = 0.00000001427746876281333
$ws.Range("S10").Value = 0.00000001427746876281333
$ws.Range("T10").Value = 58.76224640224068
$ws.Range("U10").Value = "[44.71009926216548, 72.81439354231588]"
$ws.Range("V10").Value = 0.00000000008436606968587057
$ws.Range("W10").Value = 0.00000000008436606968587057
$ws.Range("X10").Value = 17.95767767767785
$ws.Range("Y10").Value = 16.36352352352367
$ws.Range("Z10").Value = 19.55183183183202

# Row 11
$ws.Range("F11").Value = 23.42000000000022
$ws.Range("H11").Value = 0.0003999169896207189
$ws.Range("I11").Value = 0.0003999169896207189
$ws.Range("L11").Value = 43.90587041441554
$ws.Range("M11").Value = "[17.64768761503383, 70.16405321379725]"
$ws.Range("N11").Value = 0.001561045012891782
$ws.Range("O11").Value = 0.001561045012891782
$ws.Range("P11").Value = 1.880552959806503
$ws.Range("Q11").Value = "[1.2516054816103495, 2.5095004380026573]"
$ws.Range("R11").Value = 0.0000002894500148453716
$ws.Range("S11").Value = 0.0000002894500148453716
$ws.Range("T11").Value = 60.96098224997148
$ws.Range("U11").Value = "[46.649379893293684, 75.27258460664927]"
$ws.Range("V11").Value = 0.00000000005026112859241039
$ws.Range("W11").Value = 0.00000000005026112859241039
$ws.Range("X11").Value = 16.41041041041056
$ws.Range("Y11").Value = 14.0660660660662
$ws.Range("Z11").Value = 18.75475475475493

# Row 12
$ws.Range("F12").Value = 23.42000000000022
$ws.Range("H12").Value = 0.0001417337058668977
$ws.Range("I12").Value = 0.0001417337058668977
$ws.Range("L12").Value = 53.59993132300353
$ws.Range("M12").Value = "[24.513745575729644, 82.68611707027742]"
$ws.Range("N12").Value = 0.000564885513889557
$ws.Range("O12").Value = 0.000564885513889557
$ws.Range("P12").Value = 1.792500312859041
$ws.Range("Q12").Value = "[1.1887107337907326, 2.39628989192735]"
$ws.Range("R12").Value = 0.0000003350479220731728
$ws.Range("S12").Value = 0.0000003350479220731728
$ws.Range("T12").Value = 80.28404050999269
$ws.Range("U12").Value = "[64.09880470453757, 96.46927631544781]"
$ws.Range("V12").Value = 0.0000000000005355715870791755
$ws.Range("W12").Value = 0.0000000000005355715870791755
$ws.Range("X12").Value = 16.73861861861878
$ws.Range("Y12").Value = 14.48804804804818
$ws.Range("Z12").Value = 18.98918918918937

# Row 13
$ws.Range("F13").Value = 23.42000000000022
$ws.Range("H13").Value = 0.0000004187262073740428
$ws.Range("I13").Value = 0.0000004187262073740428
$ws.Range("L13").Value = 58.96036386473636
$ws.Range("M13").Value = "[36.42464573255184, 81.49608199692088]"
$ws.Range("N13").Value = 0.000003731364664139036
$ws.Range("O13").Value = 0.000003731364664139036
$ws.Range("P13").Value = 1.842816111114733
$ws.Range("Q13").Value = "[1.4025528763774249, 2.283079345852042]"
$ws.Range("R13").Value = 0.00000000008215184088555816
$ws.Range("S13").Value = 0.00000000008215184088555816
$ws.Range("T13").Value = 63.47492626060995
$ws.Range("U13").Value = "[50.566910700893175, 76.38294182032672]"
$ws.Range("V13").Value = 0.000000000000702327085377874
$ws.Range("W13").Value = 0.000000000000702327085377874
$ws.Range("X13").Value = 16.55107107107123
$ws.Range("Y13").Value = 14.91003003003017
$ws.Range("Z13").Value = 18.19211211211228

# Row 14
$ws.Range("F14").Value = 23.42000000000022
$ws.Range("H14").Value = 0.001945432100067035
$ws.Range("I14").Value = 0.001945432100067035
$ws.Range("L14").Value = 38.80243387782737
$ws.Range("M14").Value = "[11.353658874692961, 66.25120888096177]"
$ws.Range("N14").Value = 0.006622952130688864
$ws.Range("O14").Value = 0.006622952130688864
$ws.Range("P14").Value = 2.345974093671658
$ws.Range("Q14").Value = "[1.6918687163476567, 3.0000794709956584]"
$ws.Range("R14").Value = 0.000000004753175408822585
$ws.Range("S14").Value = 0.000000004753175408822585
$ws.Range("T14").Value = 58.14599264417458
$ws.Range("U14").Value = "[43.57915099927642, 72.71283428907275]"
$ws.Range("V14").Value = 0.0000000003018592042991486
$ws.Range("W14").Value = 0.0000000003018592042991486
$ws.Range("X14").Value = 14.67559559559573
$ws.Range("Y14").Value = 12.23747747747759
$ws.Range("Z14").Value = 17.11371371371387

Write-Output "applied 234 cell updates"